# "Setting Final Demand Total for shock"
#
# Adds the matti-OneDrive variant of the "Final Demand" folder paths in
# column E, rows 21-25, mirroring the existing carol-OneDrive paths that
# already live in column D for the same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$finalDemand          = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand"
$mergedFdProjected    = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand\Merged FD Projected"
$mergedFdHistorical   = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand\Merged FD Historical"
$totalFd              = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand\Total FD"

$ws.Range("E21").Value = $finalDemand
$ws.Range("E22").Value = $mergedFdProjected
$ws.Range("E23").Value = $mergedFdHistorical
$ws.Range("E24").Value = $finalDemand
$ws.Range("E25").Value = $totalFd

# Move the active selection to match the refreshed view of the sheet.
$ws.Range("E26").Select()
